# Fix lazy tab problem and fop access for RobiCategories.xlsx (BW Categories sheet)
# The D (WR_SR) and F (WR_YTH) columns are unified to match the E (WR_JR) column
# value for every populated Robi-factor row, clearing the few rows that used to
# carry stray / out-of-date figures, and re-flagging the two rows whose figures
# had been manually corrected (now highlighted across D:F instead of just E).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BW Categories")

$xlPasteFormats = -4122

function Copy-CellFormat($fromAddr, $toAddr) {
    $ws.Range($fromAddr).Copy()
    $ws.Range($toAddr).PasteSpecial($xlPasteFormats)
}

# --- Rows where the Robi factor (F) is no longer available: clear it and match
#     the plain (no value) look of columns D/E on the same row. ---
"F3", "F4", "F5", "F17" | ForEach-Object {
    $row = $_.Substring(1)
    Copy-CellFormat "D$row" $_
    $ws.Range($_).ClearContents()
}

# --- Rows where D and F must be brought in line with the (unchanged) E value.
#     Style for D/F stays whatever it already was (plain "4"), only the value
#     changes; F10 doesn't need a value fix since it already equalled E10. ---
$rows = @(
    @{R=6;  V=265},
    @{R=7;  V=307},
    @{R=8;  V=328},
    @{R=9;  V=351},
    @{R=11; V=395},
    @{R=12; V=397},
    @{R=13; V=393},
    @{R=15; V=436},
    @{R=18; V=179},
    @{R=19; V=206},
    @{R=20; V=213},
    @{R=21; V=246},
    @{R=24; V=259},
    @{R=27; V=332}
)
foreach ($item in $rows) {
    $r = $item.R
    $v = $item.V
    $ws.Range("D$r").Value = $v
    $ws.Range("F$r").Value = $v
}

# D10 needs fixing, E10/F10 already agree at 374
$ws.Range("D10").Value = 374

# --- Rows where F used to be empty and now must be populated, matching the
#     existing style of column E on the same row. ---
Copy-CellFormat "E14" "F14"
$ws.Range("D14").Value = 417
$ws.Range("F14").Value = 417

Copy-CellFormat "E26" "F26"
$ws.Range("D26").Value = 269
$ws.Range("F26").Value = 269

# --- Rows 22 & 23: the previously-isolated "needs review" highlight on E is now
#     applied across D:F, and all three columns are unified to the (unchanged)
#     E value. ---
Copy-CellFormat "E22" "D22"
Copy-CellFormat "E22" "F22"
$ws.Range("D22").Value = 260
$ws.Range("F22").Value = 260

Copy-CellFormat "E23" "D23"
Copy-CellFormat "E23" "F23"
$ws.Range("D23").Value = 262
$ws.Range("F23").Value = 262

# --- Row 25: the highlight that used to sit on F (flagging it as an outlier)
#     is removed now that F agrees with D/E; it reverts to the plain style. ---
Copy-CellFormat "E25" "F25"
$ws.Range("D25").Value = 260
$ws.Range("F25").Value = 260

# --- Restore the user's last on-screen selection (F2:F27, active cell F2) on
#     the "BW Categories" tab. ---
$ws.Range("F2:F27").Select()

Write-Output "RobiCategories BW Categories sheet updated"
